$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.061.59"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.789.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "226.91"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "32.29"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "2.047.36"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "11.34"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "1.777.62"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "0.625"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "34.066.98"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "67.98"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "242.43"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").Value = "161.99"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "16.20"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +3.67%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "1.398.79"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "0.653"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.04"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "79.90"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "0.919"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "13.65"
$ws.Range("E44").Value = "  +13.83%  "
$ws.Range("D45").Value = "6.12"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("E46").Value = "  +7.02%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "107.59"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "1.948.10"
$ws.Range("E51").Value = "  +0.00%  "
